# Uniswap workbook update — fills in the daily data rows that were missing
# on the "gUSD 26.06.25" sheet (rows 47-63), and extends the helper-column
# formulas (I, J, K, M, O, P, Q, R) that were only previously filled down to
# row 46 so that they now cover the newly entered rows through row 63.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("gUSD 26.06.25")

# ---------------------------------------------------------------------
# 1. Raw daily data for rows 47-63 (columns C,D,E,F,G) - previously blank.
# ---------------------------------------------------------------------
$rows = @(
    @(47, 98.978200000000001, 47.378500000000003, 11.01, 9.99, 5.54),
    @(48, 96.511200000000002, 48.132399999999997, 10.8, 8.93, 11.72),
    @(49, 95.649000000000001, 48.988799999999998, 10.8, 7.52, 10.25),
    @(50, 102.702, 49.462600000000002, 11.76, 7.68, 6.21),
    @(51, 99.729100000000003, 49.823500000000003, 11.49, 7.38, 5.26),
    @(52, 97.6233, 50.767600000000002, 11.34, 8.14, 11.38),
    @(53, 96.746300000000005, 50.938699999999997, 11.33, 7.53, 2.0099999999999998),
    @(54, 93.099100000000007, 52.011499999999998, 10.97, 8.83, 14.76),
    @(55, 92.088800000000006, 53.161700000000003, 10.94, 9.43, 16.09),
    @(56, 90.948099999999997, 54.155700000000003, 10.9, 9.9, 13.48),
    @(57, 89.752099999999999, 54.419899999999998, 10.85, 9.52, 3.42),
    @(58, 88.156999999999996, 55.257199999999997, 10.74, 10.28, 11.34),
    @(59, 86.425899999999999, 55.955500000000001, 10.62, 9.9600000000000009, 9.27),
    @(60, 85.294600000000003, 56.543300000000002, 10.57, 10.78, 7.35),
    @(61, 84.926500000000004, 58.529200000000003, 10.62, 12.67, 28.98),
    @(62, 79.314899999999994, 62.142400000000002, 9.9700000000000006, 17.75, 58.53),
    @(63, 78.558700000000002, 62.697800000000001, 9.9700000000000006, 16.86, 7.85)
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
    $ws.Cells.Item($r, 6).Value = $row[4]
    $ws.Cells.Item($r, 7).Value = $row[5]
}

# ---------------------------------------------------------------------
# 2. Extend the "days since start" helper column A down to row 63.
# ---------------------------------------------------------------------
$ws.Range("A62:A63").Formula = "=A61+1"

# ---------------------------------------------------------------------
# 3. Extend column I (share of total) down through the new rows.
# ---------------------------------------------------------------------
$ws.Range("I47:I63").Formula = '=C47/$A$3'

# ---------------------------------------------------------------------
# 4. Fill the day-over-day delta / rollup columns down through row 63.
#    (They previously stopped at row 24; rows 25-46 already had raw data
#    but were missing these formulas, and rows 47-63 are brand new.)
# ---------------------------------------------------------------------
$ws.Range("J25:J63").Formula = "=C25-C24"
$ws.Range("K25:K63").Formula = "=D25-D24"
$ws.Range("M25:M63").Formula = "=C25+D25"
$ws.Range("P25:P63").Formula = "=E25-E24"
$ws.Range("Q25:Q63").Formula = "=F25-F24"
$ws.Range("R25:R63").Formula = "=G25-G24"

# Column O switched (at row 24) to referencing column B instead of the
# prior-day I value; carry that same pattern down through row 63.
$ws.Range("O25:O63").Formula = "=(B25/I24-1)*100"

# Make sure the new O/P/Q/R cells pick up the same number formats as the
# existing rows above them (these columns had no cells at all in rows
# 25-46 previously, so there is no style to inherit otherwise).
$ws.Range("O24").Copy()
$ws.Range("O25:O63").PasteSpecial(-4122)
$ws.Range("P24").Copy()
$ws.Range("P25:P63").PasteSpecial(-4122)
$ws.Range("Q24").Copy()
$ws.Range("Q25:Q63").PasteSpecial(-4122)
$ws.Range("R24").Copy()
$ws.Range("R25:R63").PasteSpecial(-4122)

# Re-apply the formulas (PasteSpecial of formats only shouldn't disturb
# them, but make sure the values/formulas are exactly as intended).
$ws.Range("J25:J63").Formula = "=C25-C24"
$ws.Range("K25:K63").Formula = "=D25-D24"
$ws.Range("M25:M63").Formula = "=C25+D25"
$ws.Range("O25:O63").Formula = "=(B25/I24-1)*100"
$ws.Range("P25:P63").Formula = "=E25-E24"
$ws.Range("Q25:Q63").Formula = "=F25-F24"
$ws.Range("R25:R63").Formula = "=G25-G24"

# ---------------------------------------------------------------------
# 5. Row-height refresh for the wrapped header rows (autofit shift seen
#    after the sheet was re-saved).
# ---------------------------------------------------------------------
$wb.Worksheets.Item("gUSD 26.12.24").Rows.Item(3).RowHeight = 78
$ws.Rows.Item(4).RowHeight = 88.5
$wb.Worksheets.Item("mPendle 27.03.25").Rows.Item(4).RowHeight = 88.5

# ---------------------------------------------------------------------
# 6. Leave the selection where the author ended up after entering the
#    last row of data.
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Range("G64").Select()
